$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 66.47695399999999
$ws.Range("H2").Value = 199.430862
$ws.Range("I2").Value = 0.04311983106164722
$ws.Range("J2").Value = 0.04311983106164721
$ws.Range("M2").Value = 3.901443333333333
$ws.Range("N2").Value = 11.70433
$ws.Range("O2").Value = 0.3798296292767435
$ws.Range("P2").Value = 0.3798296292767435
$ws.Range("Q2").Value = 259.3560690036066
$ws.Range("R2").Value = 2334.20462103246
$ws.Range("S2").Value = 0.01637818944662127
$ws.Range("T2").Value = 0.01637818944662127
# Row 3
$ws.Range("G3").Value = 66.47695399999999
$ws.Range("H3").Value = 199.430862
$ws.Range("I3").Value = 0.04311983106164722
$ws.Range("J3").Value = 0.04311983106164721
$ws.Range("O3").Value = 0.2349860865034345
$ws.Range("P3").Value = 0.2349860865034345
$ws.Range("Q3").Value = 160.453695468996
$ws.Range("R3").Value = 1444.083259220964
$ws.Range("S3").Value = 0.01013256035186572
$ws.Range("T3").Value = 0.01013256035186571
# Row 4
$ws.Range("G4").Value = 66.47695399999999
$ws.Range("H4").Value = 199.430862
$ws.Range("I4").Value = 0.04311983106164722
$ws.Range("J4").Value = 0.04311983106164721
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.9736509999999999
$ws.Range("N4").Value = 2.920953
$ws.Range("O4").Value = 0.09479094447309601
$ws.Range("P4").Value = 0.09479094447309601
$ws.Range("Q4").Value = 64.72535273905399
$ws.Range("R4").Value = 582.528174651486
$ws.Range("S4").Value = 0.004087369511853882
$ws.Range("T4").Value = 0.004087369511853882
# Row 5
$ws.Range("G5").Value = 66.47695399999999
$ws.Range("H5").Value = 199.430862
$ws.Range("I5").Value = 0.04311983106164722
$ws.Range("J5").Value = 0.04311983106164721
$ws.Range("M5").Value = 2.301273333333333
$ws.Range("N5").Value = 6.90382
$ws.Range("O5").Value = 0.2240431866833358
$ws.Range("P5").Value = 0.2240431866833358
$ws.Range("Q5").Value = 152.9816415214266
$ws.Range("R5").Value = 1376.83477369284
$ws.Range("S5").Value = 0.009660704360298528
$ws.Range("T5").Value = 0.009660704360298526
# Row 6
$ws.Range("G6").Value = 66.47695399999999
$ws.Range("H6").Value = 199.430862
$ws.Range("I6").Value = 0.04311983106164722
$ws.Range("J6").Value = 0.04311983106164721
$ws.Range("M6").Value = 0.6815196666666666
$ws.Range("N6").Value = 2.044559
$ws.Range("O6").Value = 0.06635015306339016
$ws.Range("P6").Value = 0.06635015306339016
$ws.Range("Q6").Value = 45.30535153109533
$ws.Range("R6").Value = 407.748163779858
$ws.Range("S6").Value = 0.002861007391007819
$ws.Range("T6").Value = 0.002861007391007818
# Row 7
$ws.Range("I7").Value = 0.8830494168872806
$ws.Range("J7").Value = 0.8830494168872804
$ws.Range("M7").Value = 3.901443333333333
$ws.Range("N7").Value = 11.70433
$ws.Range("O7").Value = 0.3798296292767435
$ws.Range("P7").Value = 0.3798296292767435
$ws.Range("Q7").Value = 5311.343292889589
$ws.Range("R7").Value = 47802.08963600631
$ws.Range("S7").Value = 0.3354083326493403
$ws.Range("T7").Value = 0.3354083326493403
# Row 8
$ws.Range("I8").Value = 0.8830494168872806
$ws.Range("J8").Value = 0.8830494168872804
$ws.Range("O8").Value = 0.2349860865034345
$ws.Range("P8").Value = 0.2349860865034345
$ws.Range("S8").Value = 0.2075043266634819
$ws.Range("T8").Value = 0.2075043266634819
# Row 9
$ws.Range("I9").Value = 0.8830494168872806
$ws.Range("J9").Value = 0.8830494168872804
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.9736509999999999
$ws.Range("N9").Value = 2.920953
$ws.Range("O9").Value = 0.09479094447309601
$ws.Range("P9").Value = 0.09479094447309601
$ws.Range("Q9").Value = 1325.508091910919
$ws.Range("R9").Value = 11929.57282719827
$ws.Range("S9").Value = 0.08370508824316202
$ws.Range("T9").Value = 0.08370508824316202
# Row 10
$ws.Range("I10").Value = 0.8830494168872806
$ws.Range("J10").Value = 0.8830494168872804
$ws.Range("M10").Value = 2.301273333333333
$ws.Range("N10").Value = 6.90382
$ws.Range("O10").Value = 0.2240431866833358
$ws.Range("P10").Value = 0.2240431866833358
$ws.Range("Q10").Value = 3132.905348047859
$ws.Range("R10").Value = 28196.14813243074
$ws.Range("S10").Value = 0.1978412053582878
$ws.Range("T10").Value = 0.1978412053582878
# Row 11
$ws.Range("I11").Value = 0.8830494168872806
$ws.Range("J11").Value = 0.8830494168872804
$ws.Range("M11").Value = 0.6815196666666666
$ws.Range("N11").Value = 2.044559
$ws.Range("O11").Value = 0.06635015306339016
$ws.Range("P11").Value = 0.06635015306339016
$ws.Range("Q11").Value = 927.8066093118568
$ws.Range("R11").Value = 8350.259483806713
$ws.Range("S11").Value = 0.0585904639730085
$ws.Range("T11").Value = 0.05859046397300849
# Row 12
$ws.Range("G12").Value = 44.831112
$ws.Range("H12").Value = 134.493336
$ws.Range("I12").Value = 0.02907940059566787
$ws.Range("J12").Value = 0.02907940059566786
$ws.Range("M12").Value = 3.901443333333333
$ws.Range("N12").Value = 11.70433
$ws.Range("O12").Value = 0.3798296292767435
$ws.Range("P12").Value = 0.3798296292767435
$ws.Range("Q12").Value = 174.90604303832
$ws.Range("R12").Value = 1574.15438734488
$ws.Range("S12").Value = 0.01104521794784244
$ws.Range("T12").Value = 0.01104521794784244
# Row 13
$ws.Range("G13").Value = 44.831112
$ws.Range("H13").Value = 134.493336
$ws.Range("I13").Value = 0.02907940059566787
$ws.Range("J13").Value = 0.02907940059566786
$ws.Range("O13").Value = 0.2349860865034345
$ws.Range("P13").Value = 0.2349860865034345
$ws.Range("Q13").Value = 108.207689425488
$ws.Range("R13").Value = 973.8692048293921
$ws.Range("S13").Value = 0.006833254543841636
$ws.Range("T13").Value = 0.006833254543841634
# Row 14
$ws.Range("G14").Value = 44.831112
$ws.Range("H14").Value = 134.493336
$ws.Range("I14").Value = 0.02907940059566787
$ws.Range("J14").Value = 0.02907940059566786
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.9736509999999999
$ws.Range("N14").Value = 2.920953
$ws.Range("O14").Value = 0.09479094447309601
$ws.Range("P14").Value = 0.09479094447309601
$ws.Range("Q14").Value = 43.64985702991199
$ws.Range("R14").Value = 392.848713269208
$ws.Range("S14").Value = 0.002756463847174868
$ws.Range("T14").Value = 0.002756463847174867
# Row 15
$ws.Range("G15").Value = 44.831112
$ws.Range("H15").Value = 134.493336
$ws.Range("I15").Value = 0.02907940059566787
$ws.Range("J15").Value = 0.02907940059566786
$ws.Range("M15").Value = 2.301273333333333
$ws.Range("N15").Value = 6.90382
$ws.Range("O15").Value = 0.2240431866833358
$ws.Range("P15").Value = 0.2240431866833358
$ws.Range("Q15").Value = 103.16864254928
$ws.Range("R15").Value = 928.5177829435199
$ws.Range("S15").Value = 0.006515041576294722
$ws.Range("T15").Value = 0.00651504157629472
# Row 16
$ws.Range("G16").Value = 44.831112
$ws.Range("H16").Value = 134.493336
$ws.Range("I16").Value = 0.02907940059566787
$ws.Range("J16").Value = 0.02907940059566786
$ws.Range("M16").Value = 0.6815196666666666
$ws.Range("N16").Value = 2.044559
$ws.Range("O16").Value = 0.06635015306339016
$ws.Range("P16").Value = 0.06635015306339016
$ws.Range("Q16").Value = 30.553284506536
$ws.Range("R16").Value = 274.979560558824
$ws.Range("S16").Value = 0.001929422680514202
$ws.Range("T16").Value = 0.001929422680514202
# Row 17
$ws.Range("G17").Value = 52.83062100000001
$ws.Range("H17").Value = 158.491863
$ws.Range("I17").Value = 0.0342682285413064
$ws.Range("J17").Value = 0.03426822854130639
$ws.Range("M17").Value = 3.901443333333333
$ws.Range("N17").Value = 11.70433
$ws.Range("O17").Value = 0.3798296292767435
$ws.Range("P17").Value = 0.3798296292767435
$ws.Range("Q17").Value = 206.11567409631
$ws.Range("R17").Value = 1855.04106686679
$ws.Range("S17").Value = 0.01301608854281513
$ws.Range("T17").Value = 0.01301608854281513
# Row 18
$ws.Range("G18").Value = 52.83062100000001
$ws.Range("H18").Value = 158.491863
$ws.Range("I18").Value = 0.0342682285413064
$ws.Range("J18").Value = 0.03426822854130639
$ws.Range("O18").Value = 0.2349860865034345
$ws.Range("P18").Value = 0.2349860865034345
$ws.Range("Q18").Value = 127.515896311554
$ws.Range("R18").Value = 1147.643066803986
$ws.Range("S18").Value = 0.00805255691632689
$ws.Range("T18").Value = 0.008052556916326888
# Row 19
$ws.Range("G19").Value = 52.83062100000001
$ws.Range("H19").Value = 158.491863
$ws.Range("I19").Value = 0.0342682285413064
$ws.Range("J19").Value = 0.03426822854130639
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.9736509999999999
$ws.Range("N19").Value = 2.920953
$ws.Range("O19").Value = 0.09479094447309601
$ws.Range("P19").Value = 0.09479094447309601
$ws.Range("Q19").Value = 51.43858696727101
$ws.Range("R19").Value = 462.9472827054391
$ws.Range("S19").Value = 0.003248317748850339
$ws.Range("T19").Value = 0.003248317748850338
# Row 20
$ws.Range("G20").Value = 52.83062100000001
$ws.Range("H20").Value = 158.491863
$ws.Range("I20").Value = 0.0342682285413064
$ws.Range("J20").Value = 0.03426822854130639
$ws.Range("M20").Value = 2.301273333333333
$ws.Range("N20").Value = 6.90382
$ws.Range("O20").Value = 0.2240431866833358
$ws.Range("P20").Value = 0.2240431866833358
$ws.Range("Q20").Value = 121.57769929074
$ws.Range("R20").Value = 1094.19929361666
$ws.Range("S20").Value = 0.007677563124387124
$ws.Range("T20").Value = 0.007677563124387123
# Row 21
$ws.Range("G21").Value = 52.83062100000001
$ws.Range("H21").Value = 158.491863
$ws.Range("I21").Value = 0.0342682285413064
$ws.Range("J21").Value = 0.03426822854130639
$ws.Range("M21").Value = 0.6815196666666666
$ws.Range("N21").Value = 2.044559
$ws.Range("O21").Value = 0.06635015306339016
$ws.Range("P21").Value = 0.06635015306339016
$ws.Range("Q21").Value = 36.005107213713
$ws.Range("R21").Value = 324.045964923417
$ws.Range("S21").Value = 0.002273702208926915
$ws.Range("T21").Value = 0.002273702208926915
# Row 22
$ws.Range("G22").Value = 16.16161433333333
$ws.Range("H22").Value = 48.484843
$ws.Range("I22").Value = 0.01048312291409786
$ws.Range("J22").Value = 0.01048312291409786
$ws.Range("M22").Value = 3.901443333333333
$ws.Range("N22").Value = 11.70433
$ws.Range("O22").Value = 0.3798296292767435
$ws.Range("P22").Value = 0.3798296292767435
$ws.Range("Q22").Value = 63.05362249668777
$ws.Range("R22").Value = 567.48260247019
$ws.Range("S22").Value = 0.003981800690124327
$ws.Range("T22").Value = 0.003981800690124326
# Row 23
$ws.Range("G23").Value = 16.16161433333333
$ws.Range("H23").Value = 48.484843
$ws.Range("I23").Value = 0.01048312291409786
$ws.Range("J23").Value = 0.01048312291409786
$ws.Range("O23").Value = 0.2349860865034345
$ws.Range("P23").Value = 0.2349860865034345
$ws.Range("Q23").Value = 39.00886831439399
$ws.Range("R23").Value = 351.079814829546
$ws.Range("S23").Value = 0.002463388027918337
$ws.Range("T23").Value = 0.002463388027918337
# Row 24
$ws.Range("G24").Value = 16.16161433333333
$ws.Range("H24").Value = 48.484843
$ws.Range("I24").Value = 0.01048312291409786
$ws.Range("J24").Value = 0.01048312291409786
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 0.6666666666666666
$ws.Range("M24").Value = 0.9736509999999999
$ws.Range("N24").Value = 2.920953
$ws.Range("O24").Value = 0.09479094447309601
$ws.Range("P24").Value = 0.09479094447309601
$ws.Range("Q24").Value = 15.73577195726433
$ws.Range("R24").Value = 141.621947615379
$ws.Range("S24").Value = 0.000993705122054891
$ws.Range("T24").Value = 0.0009937051220548908
# Row 25
$ws.Range("G25").Value = 16.16161433333333
$ws.Range("H25").Value = 48.484843
$ws.Range("I25").Value = 0.01048312291409786
$ws.Range("J25").Value = 0.01048312291409786
$ws.Range("M25").Value = 2.301273333333333
$ws.Range("N25").Value = 6.90382
$ws.Range("O25").Value = 0.2240431866833358
$ws.Range("P25").Value = 0.2240431866833358
$ws.Range("Q25").Value = 37.19229208891777
$ws.Range("R25").Value = 334.73062880026
$ws.Range("S25").Value = 0.002348672264067582
$ws.Range("T25").Value = 0.002348672264067582
# Row 26
$ws.Range("G26").Value = 16.16161433333333
$ws.Range("H26").Value = 48.484843
$ws.Range("I26").Value = 0.01048312291409786
$ws.Range("J26").Value = 0.01048312291409786
$ws.Range("M26").Value = 0.6815196666666666
$ws.Range("N26").Value = 2.044559
$ws.Range("O26").Value = 0.06635015306339016
$ws.Range("P26").Value = 0.06635015306339016
$ws.Range("Q26").Value = 45.30535153109533
$ws.Range("R26").Value = 99.130122119237
$ws.Range("S26").Value = 0.000695556809932726
$ws.Range("T26").Value = 0.0006955568099327259
